# sd_covid_dataset.xlsx - "Updated with data from 3/22/2020"
#
# 1. Insert two new city columns ("lamesa" before nationalcity, "lakeside"
#    before ramona) into Sheet1, shifting the later city columns right.
# 2. Append a new data row (row 18, date 3/22/2020) with the day's counts.
# 3. Update the _FilterDatabase defined name to the new used range.
# 4. Remove the now-unused Sheet2 (its columns were folded into Sheet1
#    a while back; the data duplicated there is stale).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- 1. Insert the two new city columns -----------------------------------
# Before insert: BG = nationalcity, BP = ramona, BQ = ranchosantafe
$ws.Columns("BG:BG").Insert()      # new blank column -> becomes "lamesa"
$ws.Columns("BQ:BQ").Insert()      # new blank column -> becomes "lakeside"

$ws.Range("BG1").Value = "lamesa"
$ws.Range("BQ1").Value = "lakeside"

# --- 2. Append the new row for 3/22/2020 -----------------------------------
$ws.Cells.Item(18, 1).Value = "3/22/2020"

$rowValues = @(178,0,2,38,49,33,27,10,15,3,1,61,116,1,32,14,1,11,0,0,0,0,0,0,4,2,5,0,4,7,0,6,0,0,16,0,1,6,3,2,0,1,3,0,0,6,10,0,3,1,0,12,4,4,2,5,3,1,2,4,2,118,3,1,2,2,2,1,2,3)

$col = 2
foreach ($v in $rowValues) {
    $ws.Cells.Item(18, $col).Value = $v
    $col = $col + 1
}

# --- 3. Update the filter-database defined name to the new extent ----------
$wb.Names.Item("Sheet1!_FilterDatabase").RefersTo = "=Sheet1!`$A`$1:`$BS`$17"

# --- 4. Select the cell the author last had active, then drop Sheet2 -------
[void]$ws.Range("BL8").Select()

$excel.DisplayAlerts = $false
[void]$wb.Worksheets.Item("Sheet2").Delete()
$excel.DisplayAlerts = $true
